$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29; this shifts the existing rows 29-37 down to 30-38
$ws.Rows.Item(29).Insert()

# Copy the (now shifted) row 30 formatting/values into the new blank row 29
# so the row matches the existing record layout (text values, styles, etc.)
$ws.Range("A30:R30").Copy()
$ws.Range("A29:R29").PasteSpecial()

# Now overwrite the new row 29 with the new weekly record's data
$ws.Range("D29").Value = 44795
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 2700
$ws.Range("L29").Value = 2700
$ws.Range("M29").Value = 2700
$ws.Range("P29").Value = 2700
